$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.001754667048134761
$ws.Range("C2").Value = 0.3375848360084654
$ws.Range("D2").Value = 116886.6739907443
$ws.Range("E2").Value = 5548678842208.939
$ws.Range("G2").Value = 5548678959095.953
